# Generate Report for Handoff
#
# The file e2e\9f945c4e-42af-42cd-90b2-35194e932408.md has moved from
# "Handed back: in sync with en-US" to "Ready for handoff" (a newer commit
# landed upstream after the handback, so the handoff report now flags the
# handback as stale). Update the three report sheets accordingly.

$wb = $excel.ActiveWorkbook

$statusReady   = "Ready for handoff"
$overviewDate  = "2016-08-16 00:44:10"
$zhHandoffDate = "2016-08-16 00:44:02"
$errorDetail   = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ab7eccb140de83cb9813273c17b4bdb41160b377/e2e/9f945c4e-42af-42cd-90b2-35194e932408.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/01e1e43ee468f39541c97c23b230834d54a1122d/e2e/9f945c4e-42af-42cd-90b2-35194e932408.md."

# --- Overview sheet: row 3 is the 9f945c4e...md file ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = $statusReady
$wsOverview.Range("F3").Value = $statusReady
$wsOverview.Range("G3").Value = $overviewDate

# Column P ("Error Detail") widens to fit the new message. The saved OOXML
# `width` attribute is ColumnWidth + 5/6 (Excel's MDW padding), so back that
# out here to land on an exact stored width of 40.
$errorColumnWidth = 40 - (5 / 6)

# --- zh-cn sheet: row 3 is the 9f945c4e...md file ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = $statusReady
$wsZhCn.Range("H3").Value = $zhHandoffDate
$wsZhCn.Range("P3").Value = $errorDetail
$wsZhCn.Columns.Item(16).ColumnWidth = $errorColumnWidth

# --- de-de sheet: row 3 is the 9f945c4e...md file ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = $statusReady
$wsDeDe.Range("H3").Value = $overviewDate
$wsDeDe.Range("P3").Value = $errorDetail
$wsDeDe.Columns.Item(16).ColumnWidth = $errorColumnWidth
